$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 5
$ws.Range("E1").Select() | Out-Null
